# Adds new login-failure test rows (wrong user / wrong password) to the
# "LoginOnly" sheet, plus a new "ExpectedResult" column.
# Mirrors commit: add ReturnsandExchanges.java test data rows.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("LoginOnly")
$ws1 = $wb.Worksheets.Item("InputAuthorizationData")

# ---------------------------------------------------------------------
# Phase 1: write new literal values in first-use order so the shared
# string table gets the same append order as the target workbook.
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "wrongUser@rightpassword.com"
$ws.Range("B3").Value = "everlast#123"
$ws.Range("D1").Value = "ExpectedResult"
$ws.Range("D2").Value = "The account sign-in was incorrect or your account is disabled temporarily."

$ws.Range("C3").Value = "https://www.everlast.com/"
$ws.Range("D3").Value = "The account sign-in was incorrect or your account is disabled temporarily."
$ws.Range("A4").Value = "blubd.softtech@gmail.com"
$ws.Range("B4").Value = "everlast#123"
$ws.Range("C4").Value = "https://www.everlast.com/"

# ---------------------------------------------------------------------
# Phase 2: hyperlinks (order fixes the relationship id sequence)
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:wrongUser@rightpassword.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.everlast.com/")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.everlast.com/")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:blubd.softtech@gmail.com")

# ---------------------------------------------------------------------
# Phase 3: formatting. Hyperlinks.Add() stamps its own built-in
# Hyperlink style onto the cell, so re-apply the intended look afterwards
# by cloning formats from equivalent existing cells.
# ---------------------------------------------------------------------

# D1 ("ExpectedResult" header) matches the existing header style used by
# Sheet1!I1 ("expectedResult").
$ws1.Range("I1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

# A3, C3, C4 reuse the same hyperlink look already used by C2.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

# A4 reuses the same hyperlink look already used by A2.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null

# B3 / B4 ("everlast#123") are plain, unstyled text.
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

# D2 / D3 (error message) use a new red-text style.
$ws.Range("D2").Font.Color = 2567136
$ws.Range("D2").HorizontalAlignment = -4131
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Phase 4: restore final selection
# ---------------------------------------------------------------------
$ws.Range("D3").Select() | Out-Null
